# Generate Report for Handback
# Updates "latest generated" timestamps across the Overview, zh-cn, and de-de
# sheets to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" (also shared with de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-30 09:37:16"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-30 09:37:10"
$wsZhCn.Range("K2").Value = "2016-08-30 09:37:35"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-30 09:37:43"
